$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a cell to retain literal text for number-like strings (preserve
# formatting such as trailing zeros / leading zeros) without relying on the
# default General auto-detection that would otherwise coerce it to a number.
function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

$ws.Range("D2").Value = '29.275.41'
$ws.Range("E2").Value = '  +0.89%  '
$ws.Range("D3").Value = '1.861.37'
$ws.Range("E3").Value = '  +0.93%  '
Set-TextValue $ws "D4" '0.9997'
$ws.Range("E4").Value = '  -0.06%  '
Set-TextValue $ws "D5" '0.7019'
$ws.Range("E5").Value = '  +0.15%  '
Set-TextValue $ws "D6" '238.18'
$ws.Range("E6").Value = '  +0.16%  '
Set-TextValue $ws "D7" '1.000'
$ws.Range("E7").Value = '  +0.01%  '
Set-TextValue $ws "D8" '0.08274'
$ws.Range("E8").Value = '  +10.98%  '
Set-TextValue $ws "D9" '0.3048'
$ws.Range("E9").Value = '  +0.57%  '
Set-TextValue $ws "D10" '23.36'
$ws.Range("E10").Value = '  +0.34%  '
Set-TextValue $ws "D11" '0.08180'
$ws.Range("E11").Value = '  +0.77%  '
$ws.Range("D12").Value = '1.854.05'
$ws.Range("E12").Value = '  +0.15%  '
Set-TextValue $ws "D13" '0.7179'
$ws.Range("E13").Value = '  -0.75%  '
$ws.Range("E14").Value = '  -0.51%  '
Set-TextValue $ws "D15" '89.31'
$ws.Range("E15").Value = '  +0.57%  '
$ws.Range("D16").Value = '29.286.72'
$ws.Range("E16").Value = '  +0.82%  '
Set-TextValue $ws "D17" '5.789'
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws "D18" '0.000007874'
$ws.Range("E18").Value = '  +2.82%  '
$ws.Range("B19").Value = 'Avalanche'
$ws.Range("C19").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws "D19" '13.41'
$ws.Range("E19").Value = '  +3.20%  '
Set-TextValue $ws "D20" '237.09'
$ws.Range("E20").Value = '  -0.63%  '
Set-TextValue $ws "D21" '0.9986'
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("D22").Value = '2.104.39'
$ws.Range("E22").Value = '  +0.93%  '
Set-TextValue $ws "D23" '1.000'
$ws.Range("E23").Value = '  -0.04%  '
Set-TextValue $ws "D24" '7.462'
$ws.Range("E24").Value = '  -0.93%  '
Set-TextValue $ws "D25" '162.17'
$ws.Range("E25").Value = '  +0.81%  '
Set-TextValue $ws "D26" '9.002'
$ws.Range("E26").Value = '  +0.82%  '
Set-TextValue $ws "D27" '0.1456'
$ws.Range("E27").Value = '  -0.10%  '
Set-TextValue $ws "D28" '18.12'
$ws.Range("E28").Value = '  +1.01%  '
Set-TextValue $ws "D29" '1.995'
$ws.Range("E29").Value = '  +3.19%  '
Set-TextValue $ws "D30" '1.439'
$ws.Range("E30").Value = '  +4.55%  '
Set-TextValue $ws "D31" '4.423'
$ws.Range("E31").Value = '  -1.64%  '
Set-TextValue $ws "D32" '1.485'
$ws.Range("E32").Value = '  -0.34%  '
Set-TextValue $ws "D33" '4.059'
$ws.Range("E33").Value = '  +1.58%  '
Set-TextValue $ws "D34" '0.05215'
$ws.Range("E34").Value = '  +1.42%  '
Set-TextValue $ws "D35" '1.174'
Set-TextValue $ws "D36" '0.7083'
$ws.Range("E36").Value = '  +0.61%  '
Set-TextValue $ws "D37" '1.000'
$ws.Range("E37").Value = '  -2.78%  '
Set-TextValue $ws "D38" '2.661'
$ws.Range("E38").Value = '  +0.66%  '
Set-TextValue $ws "D39" '0.01851'
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("E40").Value = '  +2.11%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '1.147.73'
$ws.Range("E41").Value = '  +8.23%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws "D42" '0.9251'
$ws.Range("E42").Value = '  +2.75%  '
Set-TextValue $ws "D43" '5.950'
$ws.Range("E43").Value = '  -0.31%  '
Set-TextValue $ws "D44" '0.4286'
$ws.Range("E44").Value = '  +0.42%  '
Set-TextValue $ws "D45" '70.98'
$ws.Range("E45").Value = '  +1.61%  '
Set-TextValue $ws "D46" '0.9995'
Set-TextValue $ws "D47" '103.04'
$ws.Range("E47").Value = '  +1.23%  '
Set-TextValue $ws "D48" '1.779'
$ws.Range("E48").Value = '  +1.98%  '
$ws.Range("D49").Value = '2.000.95'
$ws.Range("E49").Value = '  +1.58%  '
Set-TextValue $ws "D50" '9.202'
$ws.Range("E50").Value = '  +0.61%  '
Set-TextValue $ws "D51" '6.976'
$ws.Range("E51").Value = '  -0.80%  '
